$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.942.62"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "1.817.81"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "231.49"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "41.50"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").Value = "0.307"
$ws.Range("E9").Value = "  +5.65%  "
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "2.083.91"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "1.837.78"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "11.04"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.65"
$ws.Range("E15").Value = "  +6.25%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.657"
$ws.Range("E16").Value = "  +5.47%  "
$ws.Range("D17").Value = "34.936.88"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "69.37"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "237.72"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "11.69"
$ws.Range("E21").Value = "  +5.50%  "
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +12.78%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("D25").Value = "172.33"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "7.74"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "17.37"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "1.56"
$ws.Range("E29").Value = "  +26.84%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "3.341.55"
$ws.Range("E31").Value = "  +37.53%  "
$ws.Range("D32").Value = "0.0547"
$ws.Range("E32").Value = "  +7.29%  "
$ws.Range("D33").Value = "3.87"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "3.97"
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("D35").Value = "1.76"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "92.45"
$ws.Range("E36").Value = "  +7.73%  "
$ws.Range("D37").Value = "0.674"
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Value = "1.313.99"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").Value = "1.28"
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("D42").Value = "0.982"
$ws.Range("E42").Value = "  +4.86%  "
$ws.Range("D43").Value = "14.64"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("E47").Value = "  +6.29%  "
$ws.Range("D48").Value = "0.0510"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").Value = "1.996.54"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "100.11"
$ws.Range("E51").Value = "  -0.27%  "
